$wb = $excel.ActiveWorkbook

# ===== Sheet: Overview =====
$ws = $wb.Worksheets.Item('Overview')

# Update cell values that changed
$ws.Range('A2').Value = 'c2de2976-cfd8-49c4-966d-c086ca65102a.md'
$ws.Range('A3').Value = 'a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md'
$ws.Range('B3').Value = 'Ready for handoff'
$ws.Range('C3').Value = 'Ready for handoff'

# Rebuild hyperlinks collection so display text matches new cell values
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/3dbb410bb84c9126777939947595179e881d6449/e2e/a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md', "", "", 'c2de2976-cfd8-49c4-966d-c086ca65102a.md')
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/3dbb410bb84c9126777939947595179e881d6449/e2e/c2de2976-cfd8-49c4-966d-c086ca65102a.md', "", "", 'a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md')
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/3dbb410bb84c9126777939947595179e881d6449/.localization-config', "", "", '.localization-config')

# ===== Sheet: zh-cn =====
$ws = $wb.Worksheets.Item('zh-cn')

# Update cell values that changed
$ws.Range('A2').Value = 'c2de2976-cfd8-49c4-966d-c086ca65102a.md'
$ws.Range('C2').Value = 'c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.zh-cn.xlf'
$ws.Range('D2').Value = '2016-02-22 14:32:35'
$ws.Range('E2').Value = 'c2de2976-cfd8-49c4-966d-c086ca65102a.md'
$ws.Range('F2').Value = 'c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.zh-cn.xlf'
$ws.Range('A3').Value = 'a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md'
$ws.Range('B3').Value = 'Ready for handoff'
$ws.Range('C3').Value = 'a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.zh-cn.xlf'
$ws.Range('D3').Value = '2016-02-22 14:32:35'
$ws.Range('E3').Value = 'a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md'
$ws.Range('F3').Value = 'a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.zh-cn.xlf'

# Rebuild hyperlinks collection so display text matches new cell values
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/3dbb410bb84c9126777939947595179e881d6449/e2e/a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md', "", "", 'c2de2976-cfd8-49c4-966d-c086ca65102a.md')
$ws.Hyperlinks.Add($ws.Range('C2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d33249393d74f37006568a55dfe1315f8550f696/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.zh-cn.xlf', "", "", 'c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('E2'), 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e7568d8c9dbac9b492ad59a8c9895812c4ff3f19/e2e/a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md', "", "", 'c2de2976-cfd8-49c4-966d-c086ca65102a.md')
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/11f9a4c5ea6ab85a0223af967a740320c12cd8cf/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.zh-cn.xlf', "", "", 'c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/3dbb410bb84c9126777939947595179e881d6449/e2e/c2de2976-cfd8-49c4-966d-c086ca65102a.md', "", "", 'a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md')
$ws.Hyperlinks.Add($ws.Range('C3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d33249393d74f37006568a55dfe1315f8550f696/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.zh-cn.xlf', "", "", 'a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('E3'), 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e7568d8c9dbac9b492ad59a8c9895812c4ff3f19/e2e/c2de2976-cfd8-49c4-966d-c086ca65102a.md', "", "", 'a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md')
$ws.Hyperlinks.Add($ws.Range('F3'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/11f9a4c5ea6ab85a0223af967a740320c12cd8cf/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.zh-cn.xlf', "", "", 'a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/3dbb410bb84c9126777939947595179e881d6449/.localization-config', "", "", '.localization-config')

# ===== Sheet: de-de =====
$ws = $wb.Worksheets.Item('de-de')

# Update cell values that changed
$ws.Range('A2').Value = 'c2de2976-cfd8-49c4-966d-c086ca65102a.md'
$ws.Range('C2').Value = 'c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.de-de.xlf'
$ws.Range('D2').Value = '2016-02-22 14:32:48'
$ws.Range('E2').Value = 'c2de2976-cfd8-49c4-966d-c086ca65102a.md'
$ws.Range('F2').Value = 'c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.de-de.xlf'
$ws.Range('A3').Value = 'a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md'
$ws.Range('B3').Value = 'Ready for handoff'
$ws.Range('C3').Value = 'a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.de-de.xlf'
$ws.Range('D3').Value = '2016-02-22 14:32:48'
$ws.Range('E3').Value = 'a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md'
$ws.Range('F3').Value = 'a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.de-de.xlf'

# Rebuild hyperlinks collection so display text matches new cell values
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/3dbb410bb84c9126777939947595179e881d6449/e2e/a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md', "", "", 'c2de2976-cfd8-49c4-966d-c086ca65102a.md')
$ws.Hyperlinks.Add($ws.Range('C2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5d372b2fce2ef50ba73cd1193940022e9834320e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.de-de.xlf', "", "", 'c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('E2'), 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/7acca1f69aaea064609e9a8bc886803547e81867/e2e/a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md', "", "", 'c2de2976-cfd8-49c4-966d-c086ca65102a.md')
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/6c6e04908eb66992d7dbab48b9b74486b94d1bbd/ol-handback/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.de-de.xlf', "", "", 'c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/3dbb410bb84c9126777939947595179e881d6449/e2e/c2de2976-cfd8-49c4-966d-c086ca65102a.md', "", "", 'a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md')
$ws.Hyperlinks.Add($ws.Range('C3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5d372b2fce2ef50ba73cd1193940022e9834320e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.de-de.xlf', "", "", 'a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('E3'), 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/7acca1f69aaea064609e9a8bc886803547e81867/e2e/c2de2976-cfd8-49c4-966d-c086ca65102a.md', "", "", 'a8acdbae-212f-4577-a97a-bcb26a9fcc8c.md')
$ws.Hyperlinks.Add($ws.Range('F3'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/6c6e04908eb66992d7dbab48b9b74486b94d1bbd/ol-handback/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/c2de2976-cfd8-49c4-966d-c086ca65102a.47c54cf4816b630d26e780babd0a85c8321691a0.de-de.xlf', "", "", 'a8acdbae-212f-4577-a97a-bcb26a9fcc8c.478929fd139052d24ef2396700ddefbd3d32b6b4.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/3dbb410bb84c9126777939947595179e881d6449/.localization-config', "", "", '.localization-config')

